$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title shape: "Vision Transformers for Information Extraction?" -> "LMMs for Information Extraction?"
$titleShape = $s.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Runs(1,1).Text = "LMMs for Information Extraction?"

# --- Body shape (Content Placeholder 2)
$body = $s.Shapes.Item(2)
$tr = $body.TextFrame.TextRange

# Paragraph 1: "Business Context" / ": Information Extraction within the SSD-DU Pipeline"
#           -> "Context" / ": PI45 The Large Multimodal Model Proof of Concept (LMM PoC)"
$para1 = $tr.Paragraphs(1,1)
$para1.Runs(1,1).Text = "Context"
$para1.Runs(2,1).Text = ": PI45 The Large Multimodal Model Proof of Concept (LMM PoC)"

# Paragraph 2: "Notes" / ": Welcome everyone. Today we're exploring a critical technology decision that could transform how we process tax document substantiation."
#           -> single, non-bold run: "Information Extraction within the SSD-DU Pipeline using VLMs"
$para2 = $tr.Paragraphs(2,1)
$para2.Runs(1,1).Text = ""
$para2 = $tr.Paragraphs(2,1)
$para2.Runs(1,1).Text = "Information Extraction within the SSD-DU Pipeline using VLMs"

# Paragraph 3: "The Business Challenge" (unchanged) / ": During " + "taxtime" + ", the ATO processes...document."
#           -> "The Business Challenge" / single run ": During Tax Time, the SSD-DU Pipeline processes thousands of WRE expense claim documents daily. Taxpayers submit receipts, invoices, and statements to support their deductions, and audit officers must verify these claims by extracting key information from each document."
$para3 = $tr.Paragraphs(3,1)
$para3.Runs(4,1).Text = ""
$para3 = $tr.Paragraphs(3,1)
$para3.Runs(3,1).Text = ""
$para3 = $tr.Paragraphs(3,1)
$para3.Runs(2,1).Text = ": During Tax Time, the SSD-DU Pipeline processes thousands of WRE expense claim documents daily. Taxpayers submit receipts, invoices, and statements to support their deductions, and audit officers must verify these claims by extracting key information from each document."

# Paragraph 4: "Current Reality" (unchanged) / ": This information extraction is currently automated using " -> ": The Information Extraction is currently automated using "
# (LayoutLM run and trailing run remain unchanged)
$para4 = $tr.Paragraphs(4,1)
$para4.Runs(2,1).Text = ": The Information Extraction is currently automated using "

# Paragraph 5: "Today's Question" -> "LLM PoC Question" (trailing run unchanged)
$para5 = $tr.Paragraphs(5,1)
$para5.Runs(1,1).Text = "LLM PoC Question"
